$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FACTORIZATION_TABLE")

# The "stroke_type" factorized values (C47:C49) were renumbered from
# -1/0/1 to 1/2/3.
$ws.Range("C47").Value = 1
$ws.Range("C48").Value = 2
$ws.Range("C49").Value = 3

# Reflect the updated view state: the window was scrolled down and the
# last edited cell (C49) is the active selection.
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
$ws.Range("C49").Select()
